$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert first new row at position 163 (pushes old row163.. down by one) ---
$ws.Rows.Item(163).Insert()

$ws.Cells.Item(163, 1).Value = 3
$ws.Cells.Item(163, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(163, 3).Value = "Coquimbo"
$ws.Cells.Item(163, 4).Value = 44567
$ws.Cells.Item(163, 5).Value = 5
$ws.Cells.Item(163, 6).Value = 100112009
$ws.Cells.Item(163, 7).Value = "Acelga"
$ws.Cells.Item(163, 8).Value = "Sin especificar"
$ws.Cells.Item(163, 9).Value = "Primera"
$ws.Cells.Item(163, 10).Value = 230
$ws.Cells.Item(163, 11).Value = 2300
$ws.Cells.Item(163, 12).Value = 2500
$ws.Cells.Item(163, 13).Value = 2413
$ws.Cells.Item(163, 14).Value = "`$/docena de atados (6 kilos)"
$ws.Cells.Item(163, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(163, 16).Value = 402
$ws.Cells.Item(163, 17).Value = 6
$ws.Cells.Item(163, 18).Value = "Hortaliza"

# --- Insert second new row at position 256 (in the now-current, already-shifted sheet) ---
$ws.Rows.Item(256).Insert()

$ws.Cells.Item(256, 1).Value = 3
$ws.Cells.Item(256, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(256, 3).Value = "Coquimbo"
$ws.Cells.Item(256, 4).Value = 44568
$ws.Cells.Item(256, 5).Value = 5
$ws.Cells.Item(256, 6).Value = 100112009
$ws.Cells.Item(256, 7).Value = "Acelga"
$ws.Cells.Item(256, 8).Value = "Sin especificar"
$ws.Cells.Item(256, 9).Value = "Primera"
$ws.Cells.Item(256, 10).Value = 310
$ws.Cells.Item(256, 11).Value = 2300
$ws.Cells.Item(256, 12).Value = 2500
$ws.Cells.Item(256, 13).Value = 2397
$ws.Cells.Item(256, 14).Value = "`$/docena de atados (6 kilos)"
$ws.Cells.Item(256, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(256, 16).Value = 400
$ws.Cells.Item(256, 17).Value = 6
$ws.Cells.Item(256, 18).Value = "Hortaliza"
